# Roll the quarterly table forward by one year:
#   - "1. kv. 2015" .. "4. kv. 2022" (rows 2-33) become
#     "1. kv. 2016" .. "4. kv. 2023"
#   - row values in columns B/C are left untouched
#   - the new trailing quarter (row 33) re-uses the existing
#     " 4,1" / " 6,4" text labels, same as before the edit
# Finally move the active selection to A33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startYear = 2016
$row = 2
for ($year = $startYear; $year -le 2023; $year++) {
    for ($q = 1; $q -le 4; $q++) {
        $ws.Range("A$row").Value = "$q. kv. $year"
        $row++
    }
}

$ws.Range("A33").Select() | Out-Null
